$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: B1 holds the "DOWNTIME (ms)" label; A1 stays empty but shares the
# wrapped-text style used by every label cell.
$ws.Cells.Item(1, 2).Value = "DOWNTIME (ms)"

# Row 2: UDP SENDER measurement
$ws.Cells.Item(2, 1).Value = "UDP SENDER"
$ws.Cells.Item(2, 2).Value = 4.6

# Row 3: UDP RECEIVER measurement
$ws.Cells.Item(3, 1).Value = "UDP RECEIVER"
$ws.Cells.Item(3, 2).Value = 1.8

# Wrap text on the label cells only (column A labels + the B1 header) -- the
# numeric value cells in column B (rows 2-3) keep the default style.
$ws.Range("A1").WrapText = $true
$ws.Range("B1").WrapText = $true
$ws.Range("A2").WrapText = $true
$ws.Range("A3").WrapText = $true

# Taller rows to fit the wrapped header/labels.
$ws.Rows.Item(1).RowHeight = 34
$ws.Rows.Item(2).RowHeight = 34
$ws.Rows.Item(3).RowHeight = 34

# Widen column A slightly so the labels read comfortably.
$ws.Columns.Item(1).ColumnWidth = 10

# Leave the selection where the author left it.
[void]$ws.Range("C3").Select()
